$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 227
$ws.Range("F4").Value = 4801
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 1171
$ws.Range("F13").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 112
$ws.Range("F19").Value = 3958
$ws.Range("F20").Value = 6289
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 3966
$ws.Range("F27").Value = 0
$ws.Range("F29").Value = 14
$ws.Range("F30").Value = 0
$ws.Range("F32").Value = 531
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 285
$ws.Range("F35").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 67
$ws.Range("F43").Value = 58

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 4801
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 748
$ws.Range("F12").Value = 223
$ws.Range("F13").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 112
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 6289
$ws.Range("F24").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 14
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("F41").Value = 961
$ws.Range("F43").Value = 67
$ws.Range("F44").Value = 58
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 0
